# Refresh cryptos list (Price / Volume(1h) columns) - GitHub Actions update.
# Values that read as numeric text (e.g. "1.00", "46.41") are written with a
# leading apostrophe so Excel keeps them as text instead of coercing them to
# numbers, matching the original inline-string cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.799.88"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "3.313.60"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'581.21"
$ws.Range("E5").Value = "  +3.18%  "
$ws.Range("D6").Value = "'182.43"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "  +2.42%  "
$ws.Range("D9").Value = "3.308.56"
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "'0.578"
$ws.Range("E11").Value = "  +0.62%  "
$ws.Range("D12").Value = "'46.41"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "'0.0000272"
$ws.Range("E13").Value = "  +2.94%  "
$ws.Range("D14").Value = "'629.85"
$ws.Range("E14").Value = "  +6.51%  "
$ws.Range("D15").Value = "3.849.67"
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "'8.45"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "67.911.22"
$ws.Range("E17").Value = "  +2.78%  "
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("D19").Value = "3.317.57"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "'17.68"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'10.91"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "'0.901"
$ws.Range("E22").Value = "  +0.61%  "
$ws.Range("D23").Value = "'17.59"
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("D24").Value = "'5.06"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").Value = "'96.79"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("D26").Value = "'3.99"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'2.77"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("D28").Value = "'9.58"
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("D29").Value = "'32.42"
$ws.Range("E29").Value = "  +5.96%  "
$ws.Range("D30").Value = "'8.59"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").Value = "'6.72"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").Value = "'599.27"
$ws.Range("E32").Value = "  +6.47%  "
$ws.Range("D33").Value = "3.940.65"
$ws.Range("E33").Value = "  +3.94%  "
$ws.Range("D34").Value = "'10.95"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("E35").Value = "  +1.57%  "
$ws.Range("D36").Value = "'3.51"
$ws.Range("E36").Value = "  -5.06%  "
$ws.Range("D37").Value = "'0.997"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "'55.75"
$ws.Range("E38").Value = "  -0.45%  "

# Rows 39/40 swapped rank order: Stacks now outranks Kaspa.
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'3.26"
$ws.Range("E39").Value = "  +3.92%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "'0.128"
$ws.Range("E40").Value = "  +0.99%  "
$ws.Range("D41").Value = "'2.68"
$ws.Range("E41").Value = "  +3.96%  "
$ws.Range("D42").Value = "'32.64"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").Value = "0.0₃0684"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "'0.338"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").Value = "'0.0413"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("E48").Value = "  +0.54%  "
$ws.Range("E49").Value = "  +12.84%  "
$ws.Range("D50").Value = "'2.55"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "'130.90"
$ws.Range("E51").Value = "  +2.21%  "
